# Generate Report for Handoff
# Adds two new "Ready for handoff" entries (8fb87282-5430-4dcc-a47a-29d3a7eb8323.md
# and e40b4a13-9c06-4a37-a4fb-8752cc14af5e.md) as new rows 6 & 7 on all three
# sheets (Overview, zh-cn, de-de), extending each sheet's table/dimension.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1.xml) - columns A-G
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

$overviewRows = @(
    @{ Row = 6; A = "8fb87282-5430-4dcc-a47a-29d3a7eb8323.md"; B = "e2e\8fb87282-5430-4dcc-a47a-29d3a7eb8323.md"; G = "'2016-08-30 06:43:50" },
    @{ Row = 7; A = "e40b4a13-9c06-4a37-a4fb-8752cc14af5e.md"; B = "e2e\e40b4a13-9c06-4a37-a4fb-8752cc14af5e.md"; G = "'2016-08-30 06:43:50" }
)

foreach ($entry in $overviewRows) {
    $r = $entry.Row
    $wsOverview.Range("A$r").Value = $entry.A
    $wsOverview.Range("C$r").Value = ".md"
    $wsOverview.Range("D$r").Value = ""
    $wsOverview.Range("E$r").Value = "Ready for handoff"
    $wsOverview.Range("F$r").Value = "Ready for handoff"
    $wsOverview.Range("G$r").NumberFormat = $dateFmt
    $wsOverview.Range("G$r").Value = $entry.G

    $wsOverview.Hyperlinks.Add($wsOverview.Range("B$r"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000" + $r + "/" + $entry.A, "", "", $entry.B) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2.xml) - columns A-P
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$loZh.ListRows.Add() | Out-Null

$zhRows = @(
    @{ Row = 6; A = "8fb87282-5430-4dcc-a47a-29d3a7eb8323.md"; G = "8fb87282-5430-4dcc-a47a-29d3a7eb8323.34e6f63693334ecaef5ce4cd257a528ebba3640e.zh-cn.xlf"; H = "'2016-08-30 06:43:45" },
    @{ Row = 7; A = "e40b4a13-9c06-4a37-a4fb-8752cc14af5e.md"; G = "e40b4a13-9c06-4a37-a4fb-8752cc14af5e.4386cc02ef683cf951d3347747a598e4d0f38286.zh-cn.xlf"; H = "'2016-08-30 06:43:45" }
)

foreach ($entry in $zhRows) {
    $r = $entry.Row
    $wsZh.Range("B$r").Value = ".md"
    $wsZh.Range("C$r").Value = "Ready for handoff"
    $wsZh.Range("D$r").Value = "e2e"
    $wsZh.Range("E$r").Value = "'ht"
    $wsZh.Range("F$r").Value = "'False"
    $wsZh.Range("G$r").Value = $entry.G
    $wsZh.Range("H$r").NumberFormat = $dateFmt
    $wsZh.Range("H$r").Value = $entry.H
    $wsZh.Range("I$r").Value = ""
    $wsZh.Range("J$r").Value = ""
    $wsZh.Range("K$r").NumberFormat = $dateFmt
    $wsZh.Range("K$r").Value = "'0001-01-01 00:00:00"
    $wsZh.Range("L$r").Value = ""
    $wsZh.Range("M$r").Value = "'True"
    $wsZh.Range("N$r").Value = ""
    $wsZh.Range("O$r").Value = "'False"
    $wsZh.Range("P$r").Value = ""

    $wsZh.Range("A$r").Value = $entry.A
    $wsZh.Hyperlinks.Add($wsZh.Range("A$r"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0000000000000000000000000000000000000" + $r + "/e2e/" + $entry.A, "", "", $entry.A) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3.xml) - columns A-P
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$loDe.ListRows.Add() | Out-Null

$deRows = @(
    @{ Row = 6; A = "8fb87282-5430-4dcc-a47a-29d3a7eb8323.md"; G = "8fb87282-5430-4dcc-a47a-29d3a7eb8323.34e6f63693334ecaef5ce4cd257a528ebba3640e.de-de.xlf"; H = "'2016-08-30 06:43:50" },
    @{ Row = 7; A = "e40b4a13-9c06-4a37-a4fb-8752cc14af5e.md"; G = "e40b4a13-9c06-4a37-a4fb-8752cc14af5e.4386cc02ef683cf951d3347747a598e4d0f38286.de-de.xlf"; H = "'2016-08-30 06:43:50" }
)

foreach ($entry in $deRows) {
    $r = $entry.Row
    $wsDe.Range("B$r").Value = ".md"
    $wsDe.Range("C$r").Value = "Ready for handoff"
    $wsDe.Range("D$r").Value = "e2e"
    $wsDe.Range("E$r").Value = "'ht"
    $wsDe.Range("F$r").Value = "'False"
    $wsDe.Range("G$r").Value = $entry.G
    $wsDe.Range("H$r").NumberFormat = $dateFmt
    $wsDe.Range("H$r").Value = $entry.H
    $wsDe.Range("I$r").Value = ""
    $wsDe.Range("J$r").Value = ""
    $wsDe.Range("K$r").NumberFormat = $dateFmt
    $wsDe.Range("K$r").Value = "'0001-01-01 00:00:00"
    $wsDe.Range("L$r").Value = ""
    $wsDe.Range("M$r").Value = "'True"
    $wsDe.Range("N$r").Value = ""
    $wsDe.Range("O$r").Value = "'False"
    $wsDe.Range("P$r").Value = ""

    $wsDe.Range("A$r").Value = $entry.A
    $wsDe.Hyperlinks.Add($wsDe.Range("A$r"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0000000000000000000000000000000000000" + $r + "/e2e/" + $entry.A, "", "", $entry.A) | Out-Null
}
